$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: G4 literal value change (GUI -> CÁ NHÂN) ---
$ws.Range("G4").Value = "CÁ NHÂN"

# --- Row 5 content (feature implementation data) ---

# D5: updated description text
$d5 = "Mô tả chức năng hệ thống`nRàng buộc dữ liệu bổ sung`nLiệt kê mô tả 04 lỗi tranh chấp (ERR21~ERR24)`nTổng hợp báo cáo lần 1"
$ws.Range("D5").Value = $d5

# E5: rich text - "Database: " (bold/italic/underline) + dispute list (normal) + trailing newline (bold/italic/underline)
$e5_r1 = "Database: "
$e5_r2 = "Tranh chấp`n1. sp_CapNhatMonHoc`n2. sp_ThemMonHoc`n3. sp_LayMonHoc`n4. sp_XoaMonHoc"
$e5_r3 = "`n"
$ws.Range("E5").Value = ($e5_r1 + $e5_r2 + $e5_r3)
$e5_c1 = $ws.Range("E5").Characters(1, $e5_r1.Length)
$e5_c1.Font.Bold = $true
$e5_c1.Font.Italic = $true
$e5_c1.Font.Underline = $true
$e5_c3 = $ws.Range("E5").Characters($e5_r1.Length + $e5_r2.Length + 1, $e5_r3.Length)
$e5_c3.Font.Bold = $true
$e5_c3.Font.Italic = $true
$e5_c3.Font.Underline = $true

# F5: rich text - "Giao diện" (bold/italic/underline) + GUI feature list (normal)
$f5_r1 = "Giao diện"
$f5_r2 = "`n1. Thêm môn học`n2. Tìm kiếm môn học`n3. Cập nhật môn học"
$ws.Range("F5").Value = ($f5_r1 + $f5_r2)
$f5_c1 = $ws.Range("F5").Characters(1, $f5_r1.Length)
$f5_c1.Font.Bold = $true
$f5_c1.Font.Italic = $true
$f5_c1.Font.Underline = $true

# G5..P5: stored-procedure / module notes for the 4 CRUD functions
$ws.Range("G5").Value = "(MONHOC)`nCập nhật môn học`nLấy danh sách môn học"
$ws.Range("H5").Value = "sp_CapNhatMonHoc`nsp_LayMonHoc"
$ws.Range("I5").Value = "(MONHOC)`nLấy sanh sách môn học`nCập nhật môn học"
$ws.Range("J5").Value = "sp_CapNhatMonHoc`nsp_LayMonHoc"
$ws.Range("K5").Value = "(MONHOC)`nLấy danh sách môn học`nThêm môn học"
$ws.Range("L5").Value = "sp_LayMonHoc`nsp_ThemMonHoc"
$ws.Range("M5").Value = "Không có"
$ws.Range("N5").Value = "Không có"
$ws.Range("O5").Value = "(MONHOC)`nCập nhật môn học`nXóa môn học"
$ws.Range("P5").Value = "sp_XoaMonHoc`nsp_CapNhatMonHoc"

# R5: updated percentage
$ws.Range("R5").Value = 0.15

# Row 5 height adjustment
$ws.Rows(5).RowHeight = 200.15

# --- Sheet view: scroll/selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("L6").Select()
